# Apply the feature-comparison update: insert two new comparison rows
# ("franzosa_ControlvsCD_ConvCD" and "franzosa_ControlvsUC_ConvUC") into the
# existing table, pushing the subsequent rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "franzosa_ControlvsCD_ConvCD" row, just before the existing
#     "franzosa_ControlvsCD_Fp" row (currently row 8). ---
$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8, 1).Value = "franzosa_ControlvsCD_ConvCD"
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0.33
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.67
$ws.Cells.Item(8, 8).Value = 0.67

# --- Insert "franzosa_ControlvsUC_ConvUC" row, just before the existing
#     "franzosa_ControlvsUC_Fp" row (now shifted down to row 13 after the
#     previous insertion). ---
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = "franzosa_ControlvsUC_ConvUC"
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0.33
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.67
$ws.Cells.Item(13, 8).Value = 0.67
